# Updates the two "on-chain evidence" placeholder sheets (B1, B2) with the
# newly-submitted transaction hashes, and leaves B2 as the active/selected
# sheet with its last-used cell selections, matching the author's re-upload.

$wb = $excel.ActiveWorkbook

# --- B1 tab: new evidence hashes -------------------------------------------------
$b1 = $wb.Worksheets.Item("B1")
$b1.Range("A2").Value = "B050DAA39175DEB9C51ED6370B44A0B076A2916836BB36ED4B134D8F7E65BA62"
$b1.Range("A3").Value = "B3A83CB2BBF8AA139D9AD074C91FBE5F3CE2BEBBCA62E1737DEBE3A4876ADB66"

# --- B2 tab: new evidence hashes -------------------------------------------------
$b2 = $wb.Worksheets.Item("B2")
$b2.Range("A2").Value = "15521471A9D448CF1E19EC90BBF50EF57DA9F17DF139B6EBC3899D4A7E9EF066"
$b2.Range("A3").Value = "CB66438079FBBD24F5C4C40ACD89859D044924C4CBD2D77A9F554C409AAF1D7C"

# --- restore each sheet's last selection, then leave B2 the active tab ----------
$b1.Activate()
$b1.Range("A2").Select()

$b2.Activate()
$b2.Range("A4").Select()
